# AudioSheet.xlsx update: add "Created"/"Implemented" tracking columns with
# color-coded status cells in column A, and move the "Moved" note from A28
# into the Notes column (G28) as lowercase "moved".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells -------------------------------------------------
# (Insert these first so the shared-string table order matches: Created,
#  Implemented, moved.)
$ws.Range("A2").Value = "Created"
$ws.Range("H2").Value = "Implemented"

# --- Column A status fill (red = created, yellow = highlighted) -------
$red = 255          # RGB(255,0,0) -> BGR 0x0000FF
$yellow = 65535      # RGB(255,255,0) -> BGR 0x00FFFF

$redRows = @(3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,28,30)
$yellowRows = @(27,29)

foreach ($r in $redRows) {
    $ws.Cells.Item($r, 1).Interior.Color = $red
}
foreach ($r in $yellowRows) {
    $ws.Cells.Item($r, 1).Interior.Color = $yellow
}

# --- Row 28: the old "Moved" note in A28 is removed; its note moves to
#     G28 (lowercase) -------------------------------------------------
$ws.Range("A28").ClearContents()
$ws.Range("A28").Interior.Color = $red
$ws.Range("G28").Value = "moved"

# --- Column widths ------------------------------------------------------
# Target OOXML <col> widths: A=13.21875, G=30.44140625, H:I=12.6640625
# (engine rounds ColumnWidth to the nearest 1/6 of a character, so feed it
# the closest achievable COM width for each target)
$ws.Columns.Item(1).ColumnWidth = 12.333333333333332
$ws.Columns.Item(7).ColumnWidth = 29.666666666666668
$ws.Columns.Item(8).ColumnWidth = 11.833333333333332
$ws.Columns.Item(9).ColumnWidth = 11.833333333333332

# --- Selection / navigation ---------------------------------------------
$ws.Range("A27").Select()
